# Adds additional submission forms, improves and add logic for processing forms.
# Currently all forms are being created, but not fully filled yet.
#
# On the "Submission_Request" sheet: rename the (optional) aliquots header
# and add a new "Experiment_id" column with sample data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Submission_Request")

# G1 header: "Aliquots (optional)" -> "Sample_id"
$ws.Range("G1").Value = "Sample_id"

# New H1 header "Experiment_id", bold like the other header cells
$ws.Range("H1").Value = "Experiment_id"
$ws.Range("H1").Font.Bold = $true

# New data values for row 2
$ws.Range("G2").Value = "AS17-00144"
$ws.Range("H2").Value = "Exp_123"

# Give column H a sensible explicit width similar to the other data columns
$ws.Columns.Item(8).ColumnWidth = 13.307291666666666
